# Applies the "remover duplicatas / tipo" edit described by the commit:
#  - Drop the "Tipo " column from Table1 on the "Tabela Principal" sheet
#    (this also removes the now-unused "Tipo " shared string and shrinks
#    the table from A6:G7 to A6:F7).
#  - Delete the now-empty helper column (column I) that sat between the
#    table and the small legend box on the "Tabela Principal" sheet.
#  - Fill in the newly collected "Taylor & Francis" search-result counts
#    on the "Metodologia" sheet (rows 8-10 and 19-21, columns K and S).
#    The SUM() totals in rows 11, 22 recalculate automatically.

$wb = $excel.ActiveWorkbook

$wsMain = $wb.Worksheets.Item("Tabela Principal")
$wsMeth = $wb.Worksheets.Item("Metodologia")

# --- Tabela Principal: remove the "Tipo " column from Table1 -----------
$tbl = $wsMain.ListObjects.Item("Table1")
$tbl.ListColumns.Item("Tipo ").Delete()

# --- Tabela Principal: remove the blank helper column (column I) -------
$wsMain.Columns.Item(9).Delete()

# --- Tabela Principal: restore the selection shown in the saved file ---
$null = $wsMain.Range("A7").Select()

# --- Metodologia: enter the new Taylor & Francis counts -----------------
$wsMeth.Range("K8").Value = 362
$wsMeth.Range("K9").Value = 10608
$wsMeth.Range("K10").Value = 54

$wsMeth.Range("K19").Value = 270
$wsMeth.Range("S19").Value = 70

$wsMeth.Range("K20").Value = 6488
$wsMeth.Range("S20").Value = 20

$wsMeth.Range("K21").Value = 43
$wsMeth.Range("S21").Value = 20

# --- Metodologia: keep it the active sheet/selection as in the file ----
$null = $wsMeth.Activate()
$null = $wsMeth.Range("R21").Select()

$wb.Save()
